# Timesheet and Weekly Report
# Added in time for stomach virus over last weekend under "illness"
# wk-18.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Saturday (row 16) and Sunday (row 17) "Sick" hours for the stomach virus
# that kept James off over the weekend.
$ws.Range("D16").Value = 6
$ws.Range("D17").Value = 4

# Scroll the sheet view back up a bit and move the active cell/selection,
# mirroring where the author left off working.
$ws.Range("P10").Select()
$excel.ActiveWindow.ScrollRow = 4
